# Fixed ICDC breed all testcases
# Update the StatQuery (column C) formulas for the CasesTab, SamplesTab and
# FilesTab rows on the "startup" sheet with the corrected Cypher query, and
# move the active selection to B3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newStatQuery = @'
MATCH (p:program)<--(s:study)<-[*]-(c:case)<--(demo:demographic)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (diag:diagnosis)-->(c)
OPTIONAL MATCH (f:file)-[*]->(c)
OPTIONAL MATCH (sf:file)-->(s)
WITH DISTINCT f, sf, samp AS samp, c, demo, diag, s, p
WHERE demo.breed IN ['Flat-Coated Retriever']
RETURN  
    count(distinct p) AS Programs,
    count(distinct s) AS Studies,
    count(distinct c) AS Cases,
    count(distinct samp) AS Samples,
    count(distinct f) AS `Case Files`,
    count(distinct sf) AS `Study Files`
'@

# Trim the trailing newline introduced by the here-string terminator so the
# cell value matches the source text exactly.
$newStatQuery = $newStatQuery.TrimEnd("`r", "`n")

$ws.Range("C2").Value = $newStatQuery
$ws.Range("C3").Value = $newStatQuery
$ws.Range("C4").Value = $newStatQuery

# Move the saved selection from B2 to B3, matching the author's last cursor
# position when the file was saved.
$ws.Range("B3").Select()
